# -----------------------------------------------------------------------
# PlayerPerformance_4571.xlsx update
#  1. Insert a new "Player Info" sheet before "ODI Batting"
#  2. Append a new "ODI Batting Extra" sheet after "ODI Bowling"
#  3. On "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE and replace
#     the full scorecard URL with just the numeric match code; also drop
#     the now-empty INNING_NUMBER placeholder cells for "did not bat" rows
#  4. On "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE and replace
#     the full scorecard URL with just the numeric match code
#  5. Fill in the new "Player Info" / "ODI Batting Extra" sheets
#
# NOTE: sheet object variables captured *before* a Worksheets.Add() call
# can resolve to the wrong sheet afterwards (this COM host seems to
# re-resolve by positional index rather than identity), so every sheet
# handle used below is (re-)fetched by name AFTER all Add() calls are
# done and the final tab order is settled.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1+2: create the final sheet order first
#   Player Info | ODI Batting | ODI Bowling | ODI Batting Extra
# ------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowlingSheet)
$extra.Name = "ODI Batting Extra"

# Re-fetch every handle by name now that the tab order is final.
$playerInfo = $wb.Worksheets.Item("Player Info")
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Item("ODI Batting Extra")

# ------------------------------------------------------------------
# Step 3: "ODI Batting" - MATCH_CARD_LINK -> MATCH_CODE,
#          drop empty INNING_NUMBER placeholder cells
# ------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$lastRow = $battingSheet.UsedRange.Rows.Count
$battingSheet.Range("D2:D" + $lastRow).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $link = $battingSheet.Cells.Item($r, 4).Value2
    if ($link -match "MatchCode=(\d+)") {
        $battingSheet.Cells.Item($r, 4).Value = $matches[1]
    }
    # Drop the empty INNING_NUMBER placeholder cell (did-not-bat rows)
    $inning = $battingSheet.Cells.Item($r, 2).Value2
    if ([string]::IsNullOrEmpty($inning)) {
        $battingSheet.Cells.Item($r, 2).ClearContents()
    }
}

# ------------------------------------------------------------------
# Step 4: "ODI Bowling" - MATCH_CARD_LINK -> MATCH_CODE
# ------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$lastRowBowl = $bowlingSheet.UsedRange.Rows.Count
$bowlingSheet.Range("B2:B" + $lastRowBowl).NumberFormat = "@"

for ($r = 2; $r -le $lastRowBowl; $r++) {
    $link = $bowlingSheet.Cells.Item($r, 2).Value2
    if ($link -match "MatchCode=(\d+)") {
        $bowlingSheet.Cells.Item($r, 2).Value = $matches[1]
    }
}

# ------------------------------------------------------------------
# Step 5a: fill "Player Info"
# ------------------------------------------------------------------
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160
$playerInfo.Range("A1:D1").Borders.LineStyle = 1

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4571"
$playerInfo.Range("B2").Value = "Paththamperuma Arachchige Don Lakshan Rangika Sandakan"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Left Arm Wrist Spin (Chinaman)"

# ------------------------------------------------------------------
# Step 5b: fill "ODI Batting Extra"
# ------------------------------------------------------------------
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
$extra.Range("A1:F1").Font.Bold = $true
$extra.Range("A1:F1").HorizontalAlignment = -4108
$extra.Range("A1:F1").VerticalAlignment = -4160
$extra.Range("A1:F1").Borders.LineStyle = 1

$extra.Range("A2:A21").NumberFormat = "@"
$extra.Range("C2:E21").NumberFormat = "@"

$extraData = @(
    @("4122", "9",  "",  "",  "",      "NO"),
    @("4182", "10", "0", "0", "2.59%", "NO"),
    @("4209", "",   "",  "",  "",      "NO"),
    @("4210", "9",  "",  "",  "",      "NO"),
    @("4215", "11", "",  "",  "",      "NO"),
    @("4231", "10", "1", "0", "1.84%", "NO"),
    @("4232", "10", "0", "0", "2.01%", "NO"),
    @("4233", "",   "",  "",  "",      "NO"),
    @("4261", "10", "0", "0", "1.30%", "NO"),
    @("4376", "9",  "0", "0", "",      "NO"),
    @("4413", "10", "0", "0", "1.03%", "NO"),
    @("4414", "",   "",  "",  "",      "NO"),
    @("4417", "10", "0", "0", "",      "NO"),
    @("4449", "10", "2", "0", "6.90%", "NO"),
    @("4450", "",   "",  "",  "",      "NO"),
    @("4451", "10", "",  "",  "",      "NO"),
    @("4463", "10", "1", "0", "3.57%", "NO"),
    @("4464", "",   "",  "",  "",      "NO"),
    @("4480", "11", "",  "",  "",      "NO"),
    @("4482", "10", "0", "0", "",      "NO")
)

$r = 2
foreach ($row in $extraData) {
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $extra.Cells.Item($r, 2).Value = [int]$row[1]
    }
    $extra.Cells.Item($r, 3).Value = $row[2]
    $extra.Cells.Item($r, 4).Value = $row[3]
    $extra.Cells.Item($r, 5).Value = $row[4]
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r++
}
